$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 20 — this pushes the existing
# rows 20-36 down to 22-38 (matching the diff's row-shift + new dimension
# A1:R38), while leaving rows 1-19 untouched.
$ws.Rows("20:21").Insert()

# --- Fill the newly inserted row 20 ---
$ws.Cells.Item(20,1).Value = 7
$ws.Cells.Item(20,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(20,3).Value = "Ñuble"
$ws.Cells.Item(20,4).Value = 44874
$ws.Cells.Item(20,5).Value = 16
$ws.Cells.Item(20,6).Value = 100112037
$ws.Cells.Item(20,7).Value = "Cebollín"
$ws.Cells.Item(20,8).Value = "Sin especificar"
$ws.Cells.Item(20,9).Value = "Primera"
$ws.Cells.Item(20,10).Value = 300
$ws.Cells.Item(20,11).Value = 600
$ws.Cells.Item(20,12).Value = 700
$ws.Cells.Item(20,13).Value = 650
$ws.Cells.Item(20,14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(20,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(20,16).Value = 108
$ws.Cells.Item(20,17).Value = 6
$ws.Cells.Item(20,18).Value = "Hortaliza"

# --- Fill the newly inserted row 21 ---
$ws.Cells.Item(21,1).Value = 7
$ws.Cells.Item(21,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(21,3).Value = "Ñuble"
$ws.Cells.Item(21,4).Value = 44874
$ws.Cells.Item(21,5).Value = 16
$ws.Cells.Item(21,6).Value = 100112037
$ws.Cells.Item(21,7).Value = "Cebollín"
$ws.Cells.Item(21,8).Value = "Sin especificar"
$ws.Cells.Item(21,9).Value = "Segunda"
$ws.Cells.Item(21,10).Value = 200
$ws.Cells.Item(21,11).Value = 500
$ws.Cells.Item(21,12).Value = 500
$ws.Cells.Item(21,13).Value = 500
$ws.Cells.Item(21,14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(21,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(21,16).Value = 83
$ws.Cells.Item(21,17).Value = 6
$ws.Cells.Item(21,18).Value = "Hortaliza"
